# Add a new tire sheet "Tir_430_50R38" by duplicating the existing
# "Tir_145_70R13" sheet (same layout/conditional formatting/validations),
# then updating the tire-specific labels, per the commit:
# "Update 2p0. Convention change to support multi-axle vehicles"

$wb = $excel.ActiveWorkbook

# Duplicate the last tire sheet; Excel places the copy right after it.
$src = $wb.Worksheets.Item("Tir_145_70R13")
$src.Copy($null, $src)

# The newly created sheet is now the last one in the workbook.
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Tir_430_50R38"

# Update the tire identification / source-file labels for the new tire.
$new.Range("H3").Value = "MFSwift_430_50R38"
$new.Range("H5").Value = "which('Truck_430_50R38.tir')"

# Leave the cursor on the tire-type selector cell, matching the edited file.
$new.Range("H6").Select()
